# Generate Report for Handback
# Update the "latest" timestamps recorded in the localization-status report
# to reflect the newly generated handback report.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# Latest HO Xliff Generate Date (Overview) / Latest Handoff Datetime (de-de)
# -- these two cells share the same underlying value in the report.
$overview.Range("G2").Value = "2016-08-18 09:07:29"
$dede.Range("H2").Value = "2016-08-18 09:07:29"

# zh-cn: Latest Handoff Datetime / Latest Handback DateTime for the first row
$zhcn.Range("H2").Value = "2016-08-18 09:07:03"
$zhcn.Range("K2").Value = "2016-08-18 09:07:44"

# de-de: Latest Handback DateTime for the first row
$dede.Range("K2").Value = "2016-08-18 09:07:52"
